$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting existing rows 188-274 down to 189-275
$ws.Rows("188:188").Insert()

# Populate the new row 188 with the new data record
$ws.Cells.Item(188, 1).Value = 9
$ws.Cells.Item(188, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44489
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = 100112039
$ws.Cells.Item(188, 7).Value = "Ciboulette"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 210
$ws.Cells.Item(188, 11).Value = 800
$ws.Cells.Item(188, 12).Value = 1000
$ws.Cells.Item(188, 13).Value = 900
$ws.Cells.Item(188, 14).Value = "$/docena de atados"
$ws.Cells.Item(188, 15).Value = "Región Metropolitana"
$ws.Cells.Item(188, 16).Value = 300
$ws.Cells.Item(188, 17).Value = 3
$ws.Cells.Item(188, 18).Value = "Hortaliza"

$ws.Range("D188").NumberFormat = "YYYY-MM-DD HH:MM:SS"
